$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 14507.143
$ws.Range("I11").Value = 14507.143
$ws.Range("K11").Value = 14507.143
$ws.Range("M11").Value = -14367.143

$ws.Range("H69").Value = 3012.6428
$ws.Range("I69").Value = 3006.5
$ws.Range("J69").Value = 3013.6667
$ws.Range("K69").Value = 9019.5
$ws.Range("L69").Value = 9041.000100000001
$ws.Range("M69").Value = -8145.5
$ws.Range("N69").Value = -10789.0001

$ws.Range("H72").Value = 3012.6428
$ws.Range("I72").Value = 3006.5
$ws.Range("J72").Value = 3013.6667
$ws.Range("K72").Value = 27058.5
$ws.Range("L72").Value = 27123.0003
$ws.Range("M72").Value = -22690.5
$ws.Range("N72").Value = -35859.0003

$ws.Range("H80").Value = 1014.75
$ws.Range("I80").Value = 617.3333
$ws.Range("K80").Value = 1851.9999
$ws.Range("M80").Value = -853.9999

$ws.Range("H83").Value = 1014.75
$ws.Range("I83").Value = 617.3333
$ws.Range("K83").Value = 5555.9997
$ws.Range("M83").Value = -563.9997000000003

$ws.Range("H86").Value = 1176.9048
$ws.Range("I86").Value = 961.7222
$ws.Range("J86").Value = 2468
$ws.Range("K86").Value = 961.7222
$ws.Range("L86").Value = 2468
$ws.Range("M86").Value = 161.2778
$ws.Range("N86").Value = -4714

$ws.Range("H89").Value = 1176.9048
$ws.Range("I89").Value = 961.7222
$ws.Range("J89").Value = 2468
$ws.Range("K89").Value = 4808.611
$ws.Range("L89").Value = 12340
$ws.Range("M89").Value = 807.3890000000001
$ws.Range("N89").Value = -23572

$ws.Range("H113").Value = 3203.5715
$ws.Range("I113").Value = 3340.6667
$ws.Range("J113").Value = 3020.7778
$ws.Range("K113").Value = 3340.6667
$ws.Range("L113").Value = 3020.7778
$ws.Range("M113").Value = -86.66670000000022
$ws.Range("N113").Value = -9528.7778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1716.6923
$ws.Range("I2").Value = 1514.4
$ws.Range("J2").Value = 2391
$ws.Range("K2").Value = 1514.4
$ws.Range("L2").Value = 2391
$ws.Range("M2").Value = -1401.4
$ws.Range("N2").Value = -2617

$ws.Range("H61").Value = 28630332
$ws.Range("I61").Value = 37075370
$ws.Range("J61").Value = 128339.25
$ws.Range("K61").Value = 37075370
$ws.Range("L61").Value = 128339.25
$ws.Range("M61").Value = -37075158
$ws.Range("N61").Value = -128763.25

$ws.Range("H116").Value = 1716.6923
$ws.Range("I116").Value = 1514.4
$ws.Range("J116").Value = 2391
$ws.Range("K116").Value = 1514.4
$ws.Range("L116").Value = 2391
$ws.Range("M116").Value = 779.5999999999999
$ws.Range("N116").Value = -6979

$ws.Range("H136").Value = 28630332
$ws.Range("I136").Value = 37075370
$ws.Range("J136").Value = 128339.25
$ws.Range("K136").Value = 111226110
$ws.Range("L136").Value = 385017.75
$ws.Range("M136").Value = -111223560
$ws.Range("N136").Value = -390117.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1716.6923
$ws.Range("I3").Value = 1514.4
$ws.Range("J3").Value = 2391
$ws.Range("K3").Value = 1514.4
$ws.Range("L3").Value = 2391
$ws.Range("M3").Value = -1400.4
$ws.Range("N3").Value = -2619

$ws.Range("H134").Value = 1594.4062
$ws.Range("I134").Value = 1633.2693
$ws.Range("K134").Value = 4899.8079
$ws.Range("M134").Value = -2364.8079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 34484416
$ws.Range("I58").Value = 41668360
$ws.Range("J58").Value = 1484.8
$ws.Range("K58").Value = 41668360
$ws.Range("L58").Value = 1484.8
$ws.Range("M58").Value = -41668157
$ws.Range("N58").Value = -1890.8

$ws.Range("H107").Value = 400.24
$ws.Range("I107").Value = 383.38095
$ws.Range("J107").Value = 488.75
$ws.Range("K107").Value = 383.38095
$ws.Range("L107").Value = 488.75
$ws.Range("M107").Value = 1536.61905
$ws.Range("N107").Value = -4328.75

$ws.Range("H136").Value = 34484416
$ws.Range("I136").Value = 41668360
$ws.Range("J136").Value = 1484.8
$ws.Range("K136").Value = 125005080
$ws.Range("L136").Value = 4454.4
$ws.Range("M136").Value = -125002530
$ws.Range("N136").Value = -9554.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 3000
$ws.Range("I9").Value = 3000
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 9000
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -8776
$ws.Range("N9").ClearContents()

$ws.Range("H131").Value = 946.5833
$ws.Range("I131").Value = 459.66666
$ws.Range("J131").Value = 1016.1429
$ws.Range("K131").Value = 1378.99998
$ws.Range("L131").Value = 3048.4287
$ws.Range("M131").Value = 3661.00002
$ws.Range("N131").Value = -13128.4287

$ws.Range("H132").Value = 2403.7
$ws.Range("I132").Value = 1745.2222
$ws.Range("J132").Value = 2685.9048
$ws.Range("K132").Value = 15706.9998
$ws.Range("L132").Value = 24173.1432
$ws.Range("M132").Value = -13176.9998
$ws.Range("N132").Value = -29233.1432

$ws.Range("H138").Value = 6223.75
$ws.Range("I138").Value = 2455.5557
$ws.Range("J138").Value = 11068.571
$ws.Range("K138").Value = 7366.6671
$ws.Range("L138").Value = 33205.713
$ws.Range("M138").Value = -2226.6671
$ws.Range("N138").Value = -43485.713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 10041
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H80").Value = 3779.2
$ws.Range("J80").Value = 3848.8572
$ws.Range("L80").Value = 3848.8572
$ws.Range("N80").Value = -5844.8572

$ws.Range("H83").Value = 3779.2
$ws.Range("J83").Value = 3848.8572
$ws.Range("L83").Value = 19244.286
$ws.Range("N83").Value = -29228.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1783256.8
$ws.Range("I46").Value = 3030944.5
$ws.Range("J46").Value = 845.7143
$ws.Range("K46").Value = 3030944.5
$ws.Range("L46").Value = 845.7143
$ws.Range("M46").Value = -3030756.5
$ws.Range("N46").Value = -1221.7143

$ws.Range("H61").Value = 1959.7949
$ws.Range("J61").Value = 1787.4
$ws.Range("L61").Value = 1787.4
$ws.Range("N61").Value = -2191.4

$ws.Range("H68").Value = 1570
$ws.Range("I68").Value = 1555
$ws.Range("K68").Value = 1555
$ws.Range("M68").Value = -806

$ws.Range("H71").Value = 1570
$ws.Range("I71").Value = 1555
$ws.Range("K71").Value = 7775
$ws.Range("M71").Value = -4031

$ws.Range("H82").Value = 2628.818
$ws.Range("I82").Value = 2001
$ws.Range("J82").Value = 2768.3333
$ws.Range("K82").Value = 2001
$ws.Range("L82").Value = 2768.3333
$ws.Range("M82").Value = -1640
$ws.Range("N82").Value = -3490.3333

$ws.Range("H85").Value = 2628.818
$ws.Range("I85").Value = 2001
$ws.Range("J85").Value = 2768.3333
$ws.Range("K85").Value = 2001
$ws.Range("L85").Value = 2768.3333
$ws.Range("M85").Value = -753
$ws.Range("N85").Value = -5264.3333

$ws.Range("H113").Value = 1959.7949
$ws.Range("J113").Value = 1787.4
$ws.Range("L113").Value = 1787.4
$ws.Range("N113").Value = -6127.4

$ws.Range("H136").Value = 126979.69
$ws.Range("I136").Value = 112400.555
$ws.Range("J136").Value = 145724.28
$ws.Range("K136").Value = 337201.665
$ws.Range("L136").Value = 437172.84
$ws.Range("M136").Value = -334651.665
$ws.Range("N136").Value = -442272.84

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3265.6086
$ws.Range("I62").Value = 3025.0833
$ws.Range("J62").Value = 3528
$ws.Range("K62").Value = 3025.0833
$ws.Range("L62").Value = 3528
$ws.Range("M62").Value = -2401.0833
$ws.Range("N62").Value = -4776

$ws.Range("H65").Value = 3265.6086
$ws.Range("I65").Value = 3025.0833
$ws.Range("J65").Value = 3528
$ws.Range("K65").Value = 15125.4165
$ws.Range("L65").Value = 17640
$ws.Range("M65").Value = -12005.4165
$ws.Range("N65").Value = -23880

$ws.Range("H132").Value = 43984
$ws.Range("I132").Value = 32404.875
$ws.Range("J132").Value = 68686.13
$ws.Range("K132").Value = 97214.625
$ws.Range("L132").Value = 206058.39
$ws.Range("M132").Value = -94684.625
$ws.Range("N132").Value = -211118.39
